$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "Status" value for row 8
$ws.Range("F8").Value = "PASS"

# Add a new data row (row 9) continuing the login data list
$ws.Range("C9").Value = "frbfnpkfsxoupxa@gmail.com"
$ws.Range("D9").Value = "oiypuLOKPS5"
$ws.Range("E9").Value = "pass"
